$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -16.34701573582464
$ws.Cells.Item(2, 3).Value = 2.009629660589098
$ws.Cells.Item(2, 4).Value = -16.34701573582464
$ws.Cells.Item(2, 5).Value = -16.34701573582464
$ws.Cells.Item(2, 6).Value = -16.34701573582464
$ws.Cells.Item(2, 7).Value = -16.34701573582464
$ws.Cells.Item(2, 8).Value = -16.34701573582464
$ws.Cells.Item(2, 9).Value = -16.34701573582464
$ws.Cells.Item(2, 10).Value = -16.34701573582464
$ws.Cells.Item(2, 11).Value = -16.34701573582464

$ws.Cells.Item(3, 2).Value = -16.34701573582464
$ws.Cells.Item(3, 3).Value = -16.34701573582464
$ws.Cells.Item(3, 4).Value = -16.34701573582464
$ws.Cells.Item(3, 5).Value = -16.34701573582464
$ws.Cells.Item(3, 6).Value = -16.34701573582464
$ws.Cells.Item(3, 7).Value = -16.34701573582464
$ws.Cells.Item(3, 8).Value = -16.34701573582464
$ws.Cells.Item(3, 9).Value = -16.34701573582464
$ws.Cells.Item(3, 10).Value = -16.34701573582464
$ws.Cells.Item(3, 11).Value = -16.34701573582464

$ws.Cells.Item(4, 2).Value = -16.34701573582464
$ws.Cells.Item(4, 3).Value = 2.249512549111533
$ws.Cells.Item(4, 4).Value = 2.106371462875994
$ws.Cells.Item(4, 5).Value = -16.34701573582464
$ws.Cells.Item(4, 6).Value = 3.398292577095668
$ws.Cells.Item(4, 7).Value = -16.34701573582464
$ws.Cells.Item(4, 8).Value = 1.748983863777128
$ws.Cells.Item(4, 9).Value = -16.34701573582464
$ws.Cells.Item(4, 10).Value = 2.751450772665505
$ws.Cells.Item(4, 11).Value = -16.34701573582464

$ws.Cells.Item(5, 2).Value = -16.34701573582464
$ws.Cells.Item(5, 3).Value = 1.717648097809831
$ws.Cells.Item(5, 4).Value = -16.34701573582464
$ws.Cells.Item(5, 5).Value = -16.34701573582464
$ws.Cells.Item(5, 6).Value = -16.34701573582464
$ws.Cells.Item(5, 7).Value = 2.929747069805344
$ws.Cells.Item(5, 8).Value = -16.34701573582464
$ws.Cells.Item(5, 9).Value = -16.34701573582464
$ws.Cells.Item(5, 10).Value = -16.34701573582464
$ws.Cells.Item(5, 11).Value = -16.34701573582464

$ws.Cells.Item(6, 2).Value = -16.34701573582464
$ws.Cells.Item(6, 3).Value = -16.34701573582464
$ws.Cells.Item(6, 4).Value = -16.34701573582464
$ws.Cells.Item(6, 5).Value = -16.34701573582464
$ws.Cells.Item(6, 6).Value = -16.34701573582464
$ws.Cells.Item(6, 7).Value = -16.34701573582464
$ws.Cells.Item(6, 8).Value = -16.34701573582464
$ws.Cells.Item(6, 9).Value = -16.34701573582464
$ws.Cells.Item(6, 10).Value = -16.34701573582464
$ws.Cells.Item(6, 11).Value = -16.34701573582464

$ws.Cells.Item(7, 2).Value = 2.707771387383987
$ws.Cells.Item(7, 3).Value = -16.34701573582464
$ws.Cells.Item(7, 4).Value = -16.34701573582464
$ws.Cells.Item(7, 5).Value = -16.34701573582464
$ws.Cells.Item(7, 6).Value = -16.34701573582464
$ws.Cells.Item(7, 7).Value = -16.34701573582464
$ws.Cells.Item(7, 8).Value = -16.34701573582464
$ws.Cells.Item(7, 9).Value = -16.34701573582464
$ws.Cells.Item(7, 10).Value = -16.34701573582464
$ws.Cells.Item(7, 11).Value = -16.34701573582464

$ws.Cells.Item(8, 2).Value = -16.34701573582464
$ws.Cells.Item(8, 3).Value = -16.34701573582464
$ws.Cells.Item(8, 4).Value = -16.34701573582464
$ws.Cells.Item(8, 5).Value = 1.869545416788898
$ws.Cells.Item(8, 6).Value = -16.34701573582464
$ws.Cells.Item(8, 7).Value = -16.34701573582464
$ws.Cells.Item(8, 8).Value = -16.34701573582464
$ws.Cells.Item(8, 9).Value = -16.34701573582464
$ws.Cells.Item(8, 10).Value = -16.34701573582464
$ws.Cells.Item(8, 11).Value = -16.34701573582464

$ws.Cells.Item(9, 2).Value = 3.751320972823597
$ws.Cells.Item(9, 3).Value = -16.34701573582464
$ws.Cells.Item(9, 4).Value = -16.34701573582464
$ws.Cells.Item(9, 5).Value = -16.34701573582464
$ws.Cells.Item(9, 6).Value = -16.34701573582464
$ws.Cells.Item(9, 7).Value = -16.34701573582464
$ws.Cells.Item(9, 8).Value = -16.34701573582464
$ws.Cells.Item(9, 9).Value = -16.34701573582464
$ws.Cells.Item(9, 10).Value = -16.34701573582464
$ws.Cells.Item(9, 11).Value = -16.34701573582464

$ws.Cells.Item(10, 2).Value = -16.34701573582464
$ws.Cells.Item(10, 3).Value = -16.34701573582464
$ws.Cells.Item(10, 4).Value = -16.34701573582464
$ws.Cells.Item(10, 5).Value = -16.34701573582464
$ws.Cells.Item(10, 6).Value = -16.34701573582464
$ws.Cells.Item(10, 7).Value = -16.34701573582464
$ws.Cells.Item(10, 8).Value = -16.34701573582464
$ws.Cells.Item(10, 9).Value = -16.34701573582464
$ws.Cells.Item(10, 10).Value = -16.34701573582464
$ws.Cells.Item(10, 11).Value = 2.024935280335591

$ws.Cells.Item(11, 2).Value = -16.34701573582464
$ws.Cells.Item(11, 3).Value = -16.34701573582464
$ws.Cells.Item(11, 4).Value = -16.34701573582464
$ws.Cells.Item(11, 5).Value = 2.798829362593765
$ws.Cells.Item(11, 6).Value = -16.34701573582464
$ws.Cells.Item(11, 7).Value = 2.678980799431003
$ws.Cells.Item(11, 8).Value = -16.34701573582464
$ws.Cells.Item(11, 9).Value = -16.34701573582464
$ws.Cells.Item(11, 10).Value = -16.34701573582464
$ws.Cells.Item(11, 11).Value = 1.347104058772876

$ws.Cells.Item(12, 2).Value = -16.34701573582464
$ws.Cells.Item(12, 3).Value = -16.34701573582464
$ws.Cells.Item(12, 4).Value = -16.34701573582464
$ws.Cells.Item(12, 5).Value = -16.34701573582464
$ws.Cells.Item(12, 6).Value = -16.34701573582464
$ws.Cells.Item(12, 7).Value = -16.34701573582464
$ws.Cells.Item(12, 8).Value = -16.34701573582464
$ws.Cells.Item(12, 9).Value = -16.34701573582464
$ws.Cells.Item(12, 10).Value = -16.34701573582464
$ws.Cells.Item(12, 11).Value = -16.34701573582464

$ws.Cells.Item(13, 2).Value = -16.34701573582464
$ws.Cells.Item(13, 3).Value = -16.34701573582464
$ws.Cells.Item(13, 4).Value = -16.34701573582464
$ws.Cells.Item(13, 5).Value = 2.351583811982223
$ws.Cells.Item(13, 6).Value = -16.34701573582464
$ws.Cells.Item(13, 7).Value = -16.34701573582464
$ws.Cells.Item(13, 8).Value = -16.34701573582464
$ws.Cells.Item(13, 9).Value = -16.34701573582464
$ws.Cells.Item(13, 10).Value = 1.928749587126734
$ws.Cells.Item(13, 11).Value = 2.121087215041674

$ws.Cells.Item(14, 2).Value = -16.34701573582464
$ws.Cells.Item(14, 3).Value = -16.34701573582464
$ws.Cells.Item(14, 4).Value = 1.411890608583341
$ws.Cells.Item(14, 5).Value = -16.34701573582464
$ws.Cells.Item(14, 6).Value = -16.34701573582464
$ws.Cells.Item(14, 7).Value = -16.34701573582464
$ws.Cells.Item(14, 8).Value = -16.34701573582464
$ws.Cells.Item(14, 9).Value = -16.34701573582464
$ws.Cells.Item(14, 10).Value = -16.34701573582464
$ws.Cells.Item(14, 11).Value = 2.3037310158577

$ws.Cells.Item(15, 2).Value = -16.34701573582464
$ws.Cells.Item(15, 3).Value = -16.34701573582464
$ws.Cells.Item(15, 4).Value = 1.279983707617997
$ws.Cells.Item(15, 5).Value = -16.34701573582464
$ws.Cells.Item(15, 6).Value = -16.34701573582464
$ws.Cells.Item(15, 7).Value = -16.34701573582464
$ws.Cells.Item(15, 8).Value = -16.34701573582464
$ws.Cells.Item(15, 9).Value = -16.34701573582464
$ws.Cells.Item(15, 10).Value = -16.34701573582464
$ws.Cells.Item(15, 11).Value = -16.34701573582464

$ws.Cells.Item(16, 2).Value = -16.34701573582464
$ws.Cells.Item(16, 3).Value = -16.34701573582464
$ws.Cells.Item(16, 4).Value = -16.34701573582464
$ws.Cells.Item(16, 5).Value = -16.34701573582464
$ws.Cells.Item(16, 6).Value = -16.34701573582464
$ws.Cells.Item(16, 7).Value = -16.34701573582464
$ws.Cells.Item(16, 8).Value = -16.34701573582464
$ws.Cells.Item(16, 9).Value = -16.34701573582464
$ws.Cells.Item(16, 10).Value = 2.029690475019074
$ws.Cells.Item(16, 11).Value = -16.34701573582464

$ws.Cells.Item(17, 2).Value = -16.34701573582464
$ws.Cells.Item(17, 3).Value = 1.892848463806636
$ws.Cells.Item(17, 4).Value = 2.185130887305624
$ws.Cells.Item(17, 5).Value = -16.34701573582464
$ws.Cells.Item(17, 6).Value = -16.34701573582464
$ws.Cells.Item(17, 7).Value = -16.34701573582464
$ws.Cells.Item(17, 8).Value = 1.417862332514731
$ws.Cells.Item(17, 9).Value = -16.34701573582464
$ws.Cells.Item(17, 10).Value = 1.457978636898982
$ws.Cells.Item(17, 11).Value = -16.34701573582464

$ws.Cells.Item(18, 2).Value = -16.34701573582464
$ws.Cells.Item(18, 3).Value = -16.34701573582464
$ws.Cells.Item(18, 4).Value = -16.34701573582464
$ws.Cells.Item(18, 5).Value = -16.34701573582464
$ws.Cells.Item(18, 6).Value = -16.34701573582464
$ws.Cells.Item(18, 7).Value = -16.34701573582464
$ws.Cells.Item(18, 8).Value = 1.04738834140984
$ws.Cells.Item(18, 9).Value = -16.34701573582464
$ws.Cells.Item(18, 10).Value = 1.394089109319252
$ws.Cells.Item(18, 11).Value = -16.34701573582464

$ws.Cells.Item(19, 2).Value = -16.34701573582464
$ws.Cells.Item(19, 3).Value = -16.34701573582464
$ws.Cells.Item(19, 4).Value = 1.528067635876082
$ws.Cells.Item(19, 5).Value = -16.34701573582464
$ws.Cells.Item(19, 6).Value = -16.34701573582464
$ws.Cells.Item(19, 7).Value = -16.34701573582464
$ws.Cells.Item(19, 8).Value = 1.76546958752617
$ws.Cells.Item(19, 9).Value = -16.34701573582464
$ws.Cells.Item(19, 10).Value = -16.34701573582464
$ws.Cells.Item(19, 11).Value = -16.34701573582464

$ws.Cells.Item(20, 2).Value = -16.34701573582464
$ws.Cells.Item(20, 3).Value = 1.053738646350594
$ws.Cells.Item(20, 4).Value = 1.665788823180784
$ws.Cells.Item(20, 5).Value = -16.34701573582464
$ws.Cells.Item(20, 6).Value = 3.241261531171283
$ws.Cells.Item(20, 7).Value = -16.34701573582464
$ws.Cells.Item(20, 8).Value = 1.501906675524828
$ws.Cells.Item(20, 9).Value = 4.321911652739746
$ws.Cells.Item(20, 10).Value = -16.34701573582464
$ws.Cells.Item(20, 11).Value = 2.035105860371625

$ws.Cells.Item(21, 2).Value = -16.34701573582464
$ws.Cells.Item(21, 3).Value = 1.097101350783421
$ws.Cells.Item(21, 4).Value = -16.34701573582464
$ws.Cells.Item(21, 5).Value = 2.098633674784804
$ws.Cells.Item(21, 6).Value = -16.34701573582464
$ws.Cells.Item(21, 7).Value = 2.579178896987438
$ws.Cells.Item(21, 8).Value = 2.502819983440541
$ws.Cells.Item(21, 9).Value = -16.34701573582464
$ws.Cells.Item(21, 10).Value = -16.34701573582464
$ws.Cells.Item(21, 11).Value = -16.34701573582464

